$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 44018
$ws.Range("C2").Value = 23242
$ws.Range("E2").Value = 2721
$ws.Range("C3").Value = 15102
$ws.Range("D3").Value = 132
$ws.Range("B4").Value = 44018
$ws.Range("C4").Value = "'214061"
$ws.Range("D4").Value = "'18596"
$ws.Range("E4").Value = 33265
$ws.Range("F4").Value = 5199
$ws.Range("G4").Value = 30.16
$ws.Range("K4").Value = 110283
$ws.Range("L4").Value = 17048
$ws.Range("B6").Value = 44018
$ws.Range("C6").Value = 52155
$ws.Range("D6").Value = 653
$ws.Range("E6").Value = 10640
$ws.Range("F6").Value = 231
$ws.Range("G6").Value = 20.4
$ws.Range("H6").Value = 35.38
$ws.Range("B7").Value = 44018
$ws.Range("C7").Value = "'25469"
$ws.Range("D7").Value = "'189"
$ws.Range("E7").Value = "'657"
$ws.Range("B8").Value = 44018
$ws.Range("C8").Value = 17152
$ws.Range("D8").Value = 593
$ws.Range("E8").Value = 1650
$ws.Range("F8").Value = 22
$ws.Range("G8").Value = 14.14
$ws.Range("H8").Value = 4.14
$ws.Range("K8").Value = 11675
$ws.Range("L8").Value = 554
$ws.Range("B12").Value = 44018
$ws.Range("C12").Value = 13507
$ws.Range("D12").Value = 515
$ws.Range("E12").Value = 259
$ws.Range("G12").Value = 1.92
$ws.Range("B13").Value = 44018
$ws.Range("C13").Value = 48992
$ws.Range("D13").Value = 1051
$ws.Range("E13").Value = 5897
$ws.Range("G13").Value = 21.03
$ws.Range("K13").Value = 28046
$ws.Range("B14").Value = 44018
$ws.Range("C14").Value = 14407
$ws.Range("D14").Value = 61
$ws.Range("E14").Value = 1842
$ws.Range("G14").Value = 24.03
$ws.Range("K14").Value = 7664
$ws.Range("B15").Value = 44018
$ws.Range("C15").Value = 44375
$ws.Range("E15").Value = 15481
$ws.Range("G15").Value = 46.61
$ws.Range("K15").Value = 33213
$ws.Range("B16").Value = 44017
$ws.Range("C16").Value = 116570
$ws.Range("D16").Value = 3534
$ws.Range("E16").Value = 3170
$ws.Range("F16").Value = 368
$ws.Range("G16").Value = 4.7
$ws.Range("H16").Value = 11.21
$ws.Range("K16").Value = 67498
$ws.Range("L16").Value = 3283
$ws.Range("B17").Value = 44018
$ws.Range("C17").Value = 69904
$ws.Range("D17").Value = 3121
$ws.Range("E17").Value = 20043
$ws.Range("F17").Value = 1263
$ws.Range("G17").Value = 28.67
$ws.Range("H17").Value = 40.47
$ws.Range("K17").Value = 57246
$ws.Range("L17").Value = 69882
$ws.Range("B18").Value = 44017
$ws.Range("C18").Value = 31257
$ws.Range("D18").Value = 1114
$ws.Range("E18").Value = 15110
$ws.Range("F18").Value = 558
$ws.Range("G18").Value = 48.34
$ws.Range("H18").Value = 50.09
$ws.Range("B19").Value = 44018
$ws.Range("C19").Value = 87705
$ws.Range("D19").Value = 6754
$ws.Range("E19").Value = 11603
$ws.Range("K19").Value = 38615
$ws.Range("B20").Value = 44018
$ws.Range("C20").Value = 203376
$ws.Range("D20").Value = 3778
$ws.Range("E20").Value = 26511
$ws.Range("F20").Value = 737
$ws.Range("G20").Value = 13.04
$ws.Range("B21").Value = 44018
$ws.Range("C21").Value = 1249
$ws.Range("G21").Value = 0.48
$ws.Range("B22").Value = 44018
$ws.Range("C22").Value = 1251
$ws.Range("E22").Value = 127
$ws.Range("G22").Value = 10.44
$ws.Range("K22").Value = 1217
$ws.Range("B23").Value = 44018
$ws.Range("C23").Value = 34257
$ws.Range("D23").Value = 1691
$ws.Range("F23").Value = 110
$ws.Range("G23").Value = 6.42
$ws.Range("H23").Value = 6.76
$ws.Range("K23").Value = 28159
$ws.Range("L23").Value = 1628
$ws.Range("B24").Value = 44018
$ws.Range("C24").Value = 20046
$ws.Range("D24").Value = 283
$ws.Range("E24").Value = 1184
$ws.Range("G24").Value = 7.65
$ws.Range("K24").Value = 15470
$ws.Range("B25").Value = 44018
$ws.Range("C25").Value = 66089
$ws.Range("D25").Value = 5897
$ws.Range("E25").Value = 19987
$ws.Range("G25").Value = 30.24
$ws.Range("H25").Value = 39.99
$ws.Range("B26").Value = 44017
$ws.Range("C26").Value = 271684
$ws.Range("D26").Value = 6300
$ws.Range("E26").Value = 7693
$ws.Range("F26").Value = 575
$ws.Range("H26").Value = 9.199999999999999
$ws.Range("K26").Value = 177012
$ws.Range("L26").Value = 6227
$ws.Range("B27").Value = 44018
$ws.Range("C27").Value = 48331
$ws.Range("D27").Value = 2505
$ws.Range("E27").Value = 5798
$ws.Range("F27").Value = 363
$ws.Range("G27").Value = 12
$ws.Range("H27").Value = 14.49
$ws.Range("B28").Value = 44018
$ws.Range("C28").Value = 1166
$ws.Range("G28").Value = 1.23
$ws.Range("K28").Value = 2026
$ws.Range("B29").Value = 44018
$ws.Range("C29").Value = 32061
$ws.Range("E29").Value = 5487
$ws.Range("G29").Value = 19.04
$ws.Range("K29").Value = 28819
$ws.Range("B30").Value = 44018
$ws.Range("C30").Value = 97064
$ws.Range("D30").Value = 2878
$ws.Range("E30").Value = 26887
$ws.Range("F30").Value = 1357
$ws.Range("G30").Value = 27.7
$ws.Range("H30").Value = 47.15
$ws.Range("B33").Value = 44018
$ws.Range("C33").Value = 12293
$ws.Range("E33").Value = 3148
$ws.Range("G33").Value = 25.61
$ws.Range("B34").Value = 44018
$ws.Range("C34").Value = 3423
$ws.Range("D34").Value = 109
$ws.Range("G34").Value = 26.52
$ws.Range("K34").Value = 3028
$ws.Range("B35").Value = 44018
$ws.Range("C35").Value = 74529
$ws.Range("D35").Value = 1398
$ws.Range("E35").Value = 11900
$ws.Range("F35").Value = 446
$ws.Range("G35").Value = 23.78
$ws.Range("H35").Value = 33.06
$ws.Range("K35").Value = 50048
$ws.Range("L35").Value = 1349
$ws.Range("B36").Value = 44018
$ws.Range("C36").Value = 147865
$ws.Range("D36").Value = 7026
$ws.Range("E36").Value = 24783
$ws.Range("G36").Value = 16.76
$ws.Range("H36").Value = 27.87
$ws.Range("B37").Value = 44018
$ws.Range("C37").Value = 8052
$ws.Range("D37").Value = 94
$ws.Range("E37").Value = 120
$ws.Range("G37").Value = 1.49
$ws.Range("H37").Value = 1.06
$ws.Range("B38").Value = 44018
$ws.Range("C38").Value = 38569
$ws.Range("D38").Value = 1474
$ws.Range("E38").Value = 7928
$ws.Range("G38").Value = 20.56
$ws.Range("H38").Value = 8.48
$ws.Range("B39").Value = 44018
$ws.Range("C39").Value = 110137
$ws.Range("D39").Value = 8198
$ws.Range("E39").Value = 10370
$ws.Range("F39").Value = 671
$ws.Range("B40").Value = 44018
$ws.Range("C40").Value = 12436
$ws.Range("D40").Value = 356
$ws.Range("E40").Value = 3592
$ws.Range("G40").Value = 31.12
$ws.Range("H40").Value = 40.17
$ws.Range("K40").Value = 11541
$ws.Range("L40").Value = 356
